# Generate Report for Handoff
#
# The localization status report was regenerated, which:
#  - bumps the Priority of the "1916abb1-661c-4c00-937b-00de07b20355" file
#    (and the other files handed off together with it) from "low" to "ht"
#    on both the zh-cn and de-de status sheets;
#  - refreshes the "Latest Handoff Datetime" on the zh-cn sheet to the new
#    handoff timestamp;
#  - refreshes the shared "Latest HO Xliff Generate Date" / handoff
#    timestamp that is common to the Overview sheet and the de-de sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Priority column: "low" -> "ht" for rows 4-7 (zh-cn and de-de)
$wsZhCn.Range("E4:E7").Value = "ht"
$wsDeDe.Range("E4:E7").Value = "ht"

# zh-cn "Latest Handoff Datetime" refreshed for rows 4-7
$wsZhCn.Range("H4:H7").Value = "2016-08-29 18:34:08"

# Overview "Latest HO Xliff Generate Date" (rows 4-7) and de-de "Latest
# Handoff Datetime" (rows 4-7) share the same underlying timestamp string;
# update both so they stay in sync.
$wsOverview.Range("G4:G7").Value = "2016-08-29 18:34:14"
$wsDeDe.Range("H4:H7").Value = "2016-08-29 18:34:14"
